$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Phone_no rows: DataType VARCHAR(15) -> INT(10)
$ws.Range("C11").Value = "INT"
$ws.Range("D11").Value = 10

$ws.Range("C23").Value = "INT"
$ws.Range("D23").Value = 10

$ws.Range("C62").Value = "INT"
$ws.Range("D62").Value = 10

$ws.Range("C65").Value = "INT"
$ws.Range("D65").Value = 10

$ws.Range("C68").Value = "INT"
$ws.Range("D68").Value = 10

# Add Size for Zip rows (INT(5))
$ws.Range("D35").Value = 5
$ws.Range("D49").Value = 5

# Update the selection / view state
$ws.Range("D12").Select()
